# Use VarName for STUB and HEADING and to join metadata sheets (Close #19)
#
# The "VarName" column on the Variables_MD sheet used short internal codes
# (antal/fsted/sex/taar) that didn't match the human readable en_varName
# column. Those codes are replaced with the corresponding en_varName value
# so the same VarName can be used to join Variables_MD with Codelists_2MD.

$wb = $excel.ActiveWorkbook

$wsVars  = $wb.Worksheets.Item("Variables_MD")
$wsCodes = $wb.Worksheets.Item("Codelists_2MD")

# Variables_MD: column B holds VarName, column D holds en_varName.
# Row 2 is the measure/value row -> "value"; the others take en_varName.
$wsVars.Range("B2").Value = "value"
$wsVars.Range("B3").Value = "place of birth"
$wsVars.Range("B4").Value = "gender"
$wsVars.Range("B5").Value = "time"

# Codelists_2MD: column A holds VarName and must be kept in sync so the
# sheets can be joined on VarName.
$wsCodes.Range("A2").Value = "place of birth"
$wsCodes.Range("A3").Value = "place of birth"
$wsCodes.Range("A4").Value = "place of birth"
$wsCodes.Range("A5").Value = "gender"
$wsCodes.Range("A6").Value = "gender"
$wsCodes.Range("A7").Value = "gender"

# Variables_MD keeps a plain selection (no longer the "tab selected" sheet,
# no frozen/scrolled topLeftCell).
$null = $wsVars.Activate()
$null = $wsVars.Range("B5").Select()

# Codelists_2MD becomes the active sheet/selection.
$null = $wsCodes.Activate()
$null = $wsCodes.Range("A20").Select()
